$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text-formatted price/text columns keep their original text formatting
# (values like "231.46" or "44.087.96" would otherwise be auto-converted to numbers)

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "44.087.96"
$ws.Range("E2").Value = "  +1.26%  "
# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.272.59"
$ws.Range("E3").Value = "  -0.31%  "
# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  +0.07%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "231.46"
$ws.Range("E5").Value = "  -0.11%  "
# Row 6
$ws.Range("E6").Value = "  +0.83%  "
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "64.09"
$ws.Range("E7").Value = "  +4.22%  "
# Row 8
$ws.Range("E8").Value = "  +0.00%  "
# Row 9
$ws.Range("E9").Value = "  +9.03%  "
# Row 10
$ws.Range("E10").Value = "  +9.57%  "
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "57.38"
$ws.Range("E11").Value = "  -0.31%  "
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "27.22"
$ws.Range("E12").Value = "  +20.02%  "
# Row 13
$ws.Range("E13").Value = "  +2.35%  "
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.612.54"
$ws.Range("E14").Value = "  +0.16%  "
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.71"
$ws.Range("E15").Value = "  -0.49%  "
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.13"
$ws.Range("E16").Value = "  +7.20%  "
# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.840"
$ws.Range("E17").Value = "  +3.46%  "
# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.271.16"
$ws.Range("E18").Value = "  -0.59%  "
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "43.986.51"
$ws.Range("E19").Value = "  +1.61%  "
# Row 20
$ws.Range("E20").Value = "  +8.33%  "
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "73.85"
$ws.Range("E21").Value = "  +0.93%  "
# Row 22
$ws.Range("E22").Value = "  -1.82%  "
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "253.56"
$ws.Range("E23").Value = "  +1.23%  "
# Row 24
$ws.Range("E24").Value = "  -0.05%  "
# Row 25
$ws.Range("E25").Value = "  -4.04%  "
# Row 26
$ws.Range("B26").NumberFormat = "@"
$ws.Range("B26").Value = "WEMIXToken"
$ws.Range("C26").NumberFormat = "@"
$ws.Range("C26").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.37"
$ws.Range("E26").Value = "  +26.36%  "
# Row 27
$ws.Range("B27").NumberFormat = "@"
$ws.Range("B27").Value = "Toncoin"
$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.27"
$ws.Range("E27").Value = "  -4.69%  "
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.03"
$ws.Range("E28").Value = "  +2.34%  "
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "171.44"
$ws.Range("E29").Value = "  +0.99%  "
# Row 30
$ws.Range("E30").Value = "  -1.83%  "
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.83"
$ws.Range("E31").Value = "  +0.68%  "
# Row 32
$ws.Range("E32").Value = "  -7.08%  "
# Row 33
$ws.Range("E33").Value = "  +3.01%  "
# Row 34
$ws.Range("E34").Value = "  +6.29%  "
# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.83"
$ws.Range("E35").Value = "  +1.59%  "
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.93"
$ws.Range("E36").Value = "  -2.02%  "
# Row 37
$ws.Range("E37").Value = "  +4.77%  "
# Row 38
$ws.Range("E38").Value = "  +1.44%  "
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.32"
$ws.Range("E39").Value = "  -3.66%  "
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0259"
$ws.Range("E40").Value = "  +3.31%  "
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.000237"
$ws.Range("E41").Value = "  +7.20%  "
# Row 42
$ws.Range("E42").Value = "  +0.31%  "
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "17.62"
$ws.Range("E43").Value = "  +5.48%  "
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0980"
$ws.Range("E44").Value = "  +0.95%  "
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.53"
$ws.Range("E45").Value = "  +17.00%  "
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.27"
$ws.Range("E46").Value = "  -4.98%  "
# Row 47
$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = "TrustWalletToken"
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.21"
$ws.Range("E47").Value = "  -0.30%  "
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "98.36"
$ws.Range("E48").Value = "  +0.84%  "
# Row 49
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = "FTXToken"
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.44"
$ws.Range("E49").Value = "  +0.94%  "
# Row 50
$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = "ordi"
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "68.11"
$ws.Range("E50").Value = "  +9.39%  "
# Row 51
$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = "Maker"
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.449.07"
$ws.Range("E51").Value = "  -1.25%  "
